$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; existing rows 40-83 shift down to 41-84.
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new weekly price entry.
$ws.Range("A40").Value = 4
$ws.Range("B40").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C40").Value = "Los Lagos"
$ws.Range("D40").Value2 = 45174
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 100112012
$ws.Range("G40").Value = "Espinaca"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 35
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 12000
$ws.Range("N40").Value = "$/cuna 10 kilos"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 1200
$ws.Range("Q40").Value = 10
$ws.Range("R40").Value = "Hortaliza"
